# SimpleSheets 0.2.0 — add support for hyperlinks: new "Adding Links" sheet
# with a Github label + URL, each one hyperlinked to https://www.github.com/

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Adding Links"

$ws.Range("A1").Value = "Github"
$ws.Range("B1").Value = "https://www.github.com/"

$ws.Hyperlinks.Add($ws.Range("A1"), "https://www.github.com/")
$ws.Hyperlinks.Add($ws.Range("B1"), "https://www.github.com/")
